$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - Como
$ws.Range("C5").Value = 26

# Row 9 - Hellas Verona
$ws.Range("B9").Value = 34
$ws.Range("D9").Value = 40.8
$ws.Range("E9").Value = 26
$ws.Range("F9").Value = 286
$ws.Range("G9").Value = 2340
$ws.Range("H9").Value = 26
$ws.Range("O9").Value = 58
$ws.Range("P9").Value = 3
$ws.Range("Q9").Value = 0.6899999999999999
$ws.Range("R9").Value = 0.46
$ws.Range("S9").Value = 1.15
$ws.Range("T9").Value = 0.58
$ws.Range("U9").Value = 1.04

# Row 19 - Sassuolo
$ws.Range("B19").Value = 30
$ws.Range("E19").Value = 26
$ws.Range("F19").Value = 286
$ws.Range("G19").Value = 2340
$ws.Range("H19").Value = 26
$ws.Range("I19").Value = 31
$ws.Range("J19").Value = 23
$ws.Range("K19").Value = 54
$ws.Range("L19").Value = 29
$ws.Range("N19").Value = 4
$ws.Range("O19").Value = 57
$ws.Range("Q19").Value = 1.19
$ws.Range("R19").Value = 0.88
$ws.Range("S19").Value = 2.08
$ws.Range("T19").Value = 1.12
$ws.Range("U19").Value = 2
